$d = $word.ActiveDocument

$replacements = @(
    @("2025-02-17 Monday", "2025-02-18 Tuesday"),
    @("230×3=", "167×7="),
    @("121×3=", "850×2="),
    @("452×4=", "978×4="),
    @("130×5=", "370×8="),
    @("349×7=", "357×3="),
    @("219×3=", "836×6="),
    @("507×8=", "285×3="),
    @("811×7=", "476×6="),
    @("745×7=", "399×6="),
    @("461×3=", "701×4="),
    @("142×7=", "242×7="),
    @("299×3=", "634×6="),
    @("492×9=", "201×6="),
    @("323×3=", "119×5="),
    @("501×9=", "209×2="),
    @("878×2=", "815×2="),
    @("326×8=", "258×5="),
    @("384×4=", "237×2="),
    @("985×8=", "814×7="),
    @("659×7=", "868×4="),
    @("459×7=", "237×4="),
    @("380×5=", "587×8="),
    @("846×9=", "427×5="),
    @("696×6=", "991×3="),
    @("423×4=", "598×9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done"
